$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Setting_Al_cars")

# ---------------------------------------------------------------------------
# 1) Insert a new "Raw materials market" process row right after the
#    "Environment" row (old row 30 "Production" and everything below shifts
#    down by one row).
# ---------------------------------------------------------------------------
$ws.Range("A30").EntireRow.Insert()
$ws.Range("C31").Copy()
$ws.Range("C30").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Range("E30").Value2 = "Market"
$ws.Range("D30").Value2 = "Raw materials market"
$ws.Range("C30").Value2 = 1

# Renumber the process group numbers of the rows that shifted down
$ws.Range("C31").Value2 = 2
$ws.Range("C32").Value2 = 3
$ws.Range("C33").Value2 = 4
$ws.Range("C34").Value2 = 5
$ws.Range("C35").Value2 = 6

# ---------------------------------------------------------------------------
# 2) Append two new process rows ("Alloy sorting", "Scrap surplus") after
#    the "Shredding and sorting of mixed scrap" row, and fix up that row's
#    label/number.
# ---------------------------------------------------------------------------
$ws.Range("A37:A38").EntireRow.Insert()
$ws.Range("C36").Copy()
$ws.Range("C37:C38").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Range("C36").Value2 = 7
$ws.Range("D36").Value2 = "Shredding and sorting of mixed scrap"

$ws.Range("C37").Value2 = 8
$ws.Range("D37").Value2 = "Alloy sorting"
$ws.Range("E37").Value2 = "Industry/Transformation"

$ws.Range("C38").Value2 = 9
$ws.Range("D38").Value2 = "Scrap surplus"
$ws.Range("E38").Value2 = "Industry/Transformation"

# ---------------------------------------------------------------------------
# 3) Insert a new "Alloy_Sorting" parameter row (after "Dismantling", which
#    is now at row 51) replacing the old blank spacer row.
# ---------------------------------------------------------------------------
$ws.Range("A52").EntireRow.Insert()
$ws.Range("A51:I51").Copy()
$ws.Range("A52:I52").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Range("A52").Value2 = " "
$ws.Range("C52").Value2 = "Alloy_Sorting"
$ws.Range("D52").Value2 = "Alloy sorting rate"
$ws.Range("E52").Value2 = "Alloy_Sorting"
$ws.Range("F52").Value2 = "tr"
$ws.Range("G52").Value2 = "[0,1]"
$ws.Range("H52").Value2 = "[0]"
$ws.Range("I52").ClearContents()

# ---------------------------------------------------------------------------
# 4) Sheet view cosmetic updates
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 130
$ws.Range("B13").Select()
$ws.Range("A53:XFD53").Select()
